# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with freshly scraped values. Numeric-looking price strings are
# written with a leading apostrophe so Excel keeps them as plain text
# (matching the source data, which stores prices/volumes as text, not
# numbers) instead of silently converting them to numeric cell values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.822.08'
$ws.Range("E2").Value = '  +7.72%  '

$ws.Range("D3").Value = '1.751.99'
$ws.Range("E3").Value = '  +5.21%  '

$ws.Range("D4").Value = '''0.9961'
$ws.Range("E4").Value = '  -0.55%  '

$ws.Range("D5").Value = '''333.18'
$ws.Range("E5").Value = '  +5.17%  '

$ws.Range("D6").Value = '''0.9952'
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("D7").Value = '''0.3757'
$ws.Range("E7").Value = '  +3.39%  '

$ws.Range("D8").Value = '''49.92'
$ws.Range("E8").Value = '  +5.22%  '

$ws.Range("D9").Value = '''0.3456'
$ws.Range("E9").Value = '  +5.14%  '

$ws.Range("E10").Value = '  +5.72%  '

$ws.Range("D11").Value = '''0.07535'
$ws.Range("E11").Value = '  +5.95%  '

$ws.Range("D12").Value = '''0.9956'
$ws.Range("E12").Value = '  -0.30%  '

$ws.Range("D13").Value = '''6.508'
$ws.Range("E13").Value = '  +6.82%  '

$ws.Range("D14").Value = '''20.62'
$ws.Range("E14").Value = '  +4.27%  '

$ws.Range("D15").Value = '''7.101'
$ws.Range("E15").Value = '  +6.63%  '

$ws.Range("D16").Value = '1.740.20'
$ws.Range("E16").Value = '  +4.42%  '

$ws.Range("E17").Value = '  +4.32%  '

$ws.Range("D18").Value = '''0.06706'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").Value = '''84.39'
$ws.Range("E19").Value = '  +5.70%  '

$ws.Range("D20").Value = '''0.9947'
$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("D21").Value = '''16.87'
$ws.Range("E21").Value = '  +6.16%  '

$ws.Range("D22").Value = '''6.221'
$ws.Range("E22").Value = '  +4.32%  '

$ws.Range("D23").Value = '''13.20'
$ws.Range("E23").Value = '  +3.81%  '

$ws.Range("D24").Value = '26.747.00'
$ws.Range("E24").Value = '  +7.46%  '

$ws.Range("D25").Value = '''2.472'
$ws.Range("E25").Value = '  +1.40%  '

$ws.Range("D26").Value = '''2.541'
$ws.Range("E26").Value = '  +3.92%  '

$ws.Range("D27").Value = '''1.419'
$ws.Range("E27").Value = '  +14.06%  '

$ws.Range("D28").Value = '''153.23'
$ws.Range("E28").Value = '  +3.45%  '

$ws.Range("D29").Value = '''19.78'
$ws.Range("E29").Value = '  +5.58%  '

$ws.Range("D30").Value = '1.936.25'
$ws.Range("E30").Value = '  +4.59%  '

$ws.Range("D31").Value = '''132.52'
$ws.Range("E31").Value = '  +5.05%  '

$ws.Range("D32").Value = '''4.135'
$ws.Range("E32").Value = '  +0.17%  '

$ws.Range("D33").Value = '''6.267'
$ws.Range("E33").Value = '  +5.92%  '

$ws.Range("D34").Value = '''0.08641'
$ws.Range("E34").Value = '  +1.28%  '

$ws.Range("D35").Value = '''1.726'
$ws.Range("E35").Value = '  +4.00%  '

$ws.Range("D36").Value = '''13.20'
$ws.Range("E36").Value = '  +6.15%  '

$ws.Range("D37").Value = '''5.512'
$ws.Range("E37").Value = '  +5.16%  '

$ws.Range("D38").Value = '''0.02385'
$ws.Range("E38").Value = '  +4.61%  '

$ws.Range("D39").Value = '''0.06392'
$ws.Range("E39").Value = '  +4.79%  '

$ws.Range("D40").Value = '''0.2187'
$ws.Range("E40").Value = '  +4.60%  '

$ws.Range("D41").Value = '''8.715'
$ws.Range("E41").Value = '  +3.84%  '

$ws.Range("D42").Value = '''1.248'
$ws.Range("E42").Value = '  -3.11%  '

$ws.Range("D43").Value = '''0.6309'
$ws.Range("E43").Value = '  +5.34%  '

$ws.Range("D44").Value = '''14.56'
$ws.Range("E44").Value = '  +12.93%  '

$ws.Range("D45").Value = '''0.9948'
$ws.Range("E45").Value = '  -0.38%  '

$ws.Range("D46").Value = '''3.912'
$ws.Range("E46").Value = '  +2.27%  '

$ws.Range("D47").Value = '''0.6122'
$ws.Range("E47").Value = '  +7.72%  '

$ws.Range("E48").Value = '  +5.22%  '

$ws.Range("D49").Value = '''129.48'
$ws.Range("E49").Value = '  +2.61%  '

$ws.Range("D50").Value = '''0.07376'
$ws.Range("E50").Value = '  +4.83%  '

$ws.Range("D51").Value = '''78.21'
$ws.Range("E51").Value = '  +3.77%  '
